$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the style/format from H1 (IP)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0 and IF numeric values for data rows 2 through 67
$iVals = @(6,5,6,7,8,9,6,6,8,7,6,8,7,3,6,7,9,6,9,6,6,6,6,5,6,6,6,7,6,7,9,7,7,4,5,8,4,4,7,8,7,6,10,7,7,6,7,5,8,4,6,4,6,7,8,5,11,6,7,7,4,7,8,8,7,6)
$jVals = @(6,5,6,7,8,9,6,6,8,7,6,8,7,4,7,7,9,6,9,6,6,6,7,5,6,7,7,7,6,8,9,7,8,5,6,8,4,5,7,8,7,6,10,7,7,7,7,6,8,5,6,4,7,8,8,6,12,7,8,7,5,7,8,8,7,6)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
